# issue #5: stock data output to json file
#
# The "股票" (stock) sheet gains a new "property_category" column (with
# value "stock" for the existing data row), inserted just before the
# existing "date" column. This pushes the existing date/legislator_name/
# legislator_id columns one column to the right. Also fixes a stray space
# typo in the company name in the "owner"/name cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column at H (before the existing "date" column), shifting
# date / legislator_name / legislator_id from H/I/J to I/J/K.
$ws.Range("H1").EntireColumn.Insert()

# New column header + value for the stock property category.
$ws.Range("H1").Value = "property_category"
$ws.Range("H2").Value = "stock"

# Fix typo: remove stray space in the company name.
$ws.Range("B2").Value = "中日國際企業股份有限公司"
